# Sprint Review / Design Pattern Updates
# Applies the Sprint 6 content updates to the document.

$d = $word.ActiveDocument

# 1) "Features implemented:" bullet (paragraph 4) - replace "None yet"
#    with the description of the login/Firebase work in progress.
$r1 = $d.Paragraphs(4).Range
$r1.Find.Execute("None yet", $false, $false, $false, $false, $false, `
                  $true, 1, $false, `
                  "Attempts were made and are ongoing in order to create a working login system and a successful link with Firebase, which we have chosen as our database. The foundations of several other features have been written.", `
                  2) | Out-Null

# 2) "Issues fixed:" bullet (paragraph 6) - replace "None yet" with the
#    firebase fix note.
$r2 = $d.Paragraphs(6).Range
$r2.Find.Execute("None yet", $false, $false, $false, $false, $false, `
                  $true, 1, $false, `
                  "Fixes with firebase implementation are ongoing.", `
                  2) | Out-Null

# 3) Closing paragraph (last paragraph) - replace placeholder description text.
$r3 = $d.Paragraphs($d.Paragraphs.Count).Range
$r3.Find.Execute("Description Here", $false, $false, $false, $false, $false, `
                  $true, 1, $false, `
                  "Overall, we have learned over the course of this project some things which we would certainly reconsider given a second chance. Our app continues to only run at a very basic level, but hypothetical focus/usability study group scenarios have been created for when app function increases. Documentation has also been updated from throughout the project.", `
                  2) | Out-Null

# 4) Move the _GoBack bookmark from the "Issues fixed:" paragraph to the end
#    of the last paragraph (now containing the new closing text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$lastParaRange = $d.Paragraphs($d.Paragraphs.Count).Range
$d.Bookmarks.Add("_GoBack", $lastParaRange)
